# Apply the edit described by the diff:
#  - The shared string "CAA40318A" (in cell D13) is corrected to "CAA40318.1"
#  - The selected cell/range in the sheet view moves to D13
#
# All other apparent cell changes in the diff are simply a side effect of the
# shared-strings table being re-ordered after removing the old "CAA40318A"
# entry and appending the corrected "CAA40318.1" entry at the end; the actual
# displayed values in every other cell are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the erroneous accession number.
$ws.Range("D13").Value = "CAA40318.1"

# Update the selected cell/range to match the saved view state.
$ws.Range("D13").Select()

$wb.Save()
